# Updates the cryptocurrency price list (columns D = Price, E = Volume(1h))
# and, for two rows, also updates the Coin name / Link (columns B, C) to
# reflect a re-ranking between mCoin and PaxDollar.
#
# All of these columns hold plain text in the worksheet (e.g. "25.926.00",
# "0.5438", "  -0.22%  "). Excel's Range.Value setter auto-detects strings
# that look like plain numbers (one decimal point, no thousands-style extra
# dots) and silently coerces them to a real number - which would both change
# the cell's stored type and normalize/round the text (e.g. "0.8160" -> 0.816).
# To keep those cells as text - matching the original file - we force the
# cell to Text format before assigning, then restore the Normal style so we
# don't leave a stray number-format behind on cells that should look exactly
# as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)

    $looksNumeric = $Value -match '^[+-]?\d+(\.\d+)?$'

    if ($looksNumeric) {
        $Range.NumberFormat = "@"
        $Range.Value = $Value
        $Range.Style = "Normal"
    } else {
        $Range.Value = $Value
    }
}

$updates = @(
    @{ Cell = 'D2'; Value = '25.928.95' },
    @{ Cell = 'E2'; Value = '  -0.16%  ' },
    @{ Cell = 'D3'; Value = '1.639.05' },
    @{ Cell = 'E4'; Value = '  +0.29%  ' },
    @{ Cell = 'D5'; Value = '214.79' },
    @{ Cell = 'E5'; Value = '  -0.34%  ' },
    @{ Cell = 'D6'; Value = '0.5059' },
    @{ Cell = 'E6'; Value = '  -0.06%  ' },
    @{ Cell = 'D7'; Value = '1.003' },
    @{ Cell = 'E7'; Value = '  +0.17%  ' },
    @{ Cell = 'E8'; Value = '  -1.00%  ' },
    @{ Cell = 'D9'; Value = '0.06362' },
    @{ Cell = 'E9'; Value = '  -0.94%  ' },
    @{ Cell = 'E10'; Value = '  -1.01%  ' },
    @{ Cell = 'D11'; Value = '0.07734' },
    @{ Cell = 'E11'; Value = '  -0.29%  ' },
    @{ Cell = 'D12'; Value = '4.277' },
    @{ Cell = 'E12'; Value = '  +0.09%  ' },
    @{ Cell = 'D13'; Value = '1.646.90' },
    @{ Cell = 'E13'; Value = '  +0.12%  ' },
    @{ Cell = 'D14'; Value = '0.5438' },
    @{ Cell = 'E14'; Value = '  -0.62%  ' },
    @{ Cell = 'E15'; Value = '  -1.79%  ' },
    @{ Cell = 'D16'; Value = '64.08' },
    @{ Cell = 'E16'; Value = '  -0.65%  ' },
    @{ Cell = 'D17'; Value = '25.960.03' },
    @{ Cell = 'E17'; Value = '  -0.13%  ' },
    @{ Cell = 'E18'; Value = '  +0.17%  ' },
    @{ Cell = 'D19'; Value = '197.16' },
    @{ Cell = 'E19'; Value = '  -2.85%  ' },
    @{ Cell = 'D20'; Value = '4.455' },
    @{ Cell = 'E20'; Value = '  +1.37%  ' },
    @{ Cell = 'D21'; Value = '9.934' },
    @{ Cell = 'E21'; Value = '  +0.21%  ' },
    @{ Cell = 'D22'; Value = '6.018' },
    @{ Cell = 'E22'; Value = '  +0.26%  ' },
    @{ Cell = 'E23'; Value = '  +0.21%  ' },
    @{ Cell = 'D24'; Value = '1.883' },
    @{ Cell = 'E24'; Value = '  +0.24%  ' },
    @{ Cell = 'D25'; Value = '140.68' },
    @{ Cell = 'E25'; Value = '  -0.30%  ' },
    @{ Cell = 'D26'; Value = '0.1181' },
    @{ Cell = 'E26'; Value = '  +3.62%  ' },
    @{ Cell = 'E27'; Value = '  +0.60%  ' },
    @{ Cell = 'D29'; Value = '1.235' },
    @{ Cell = 'E29'; Value = '  -0.48%  ' },
    @{ Cell = 'D30'; Value = '0.04925' },
    @{ Cell = 'E30'; Value = '  -0.11%  ' },
    @{ Cell = 'D31'; Value = '3.251' },
    @{ Cell = 'E31'; Value = '  -0.71%  ' },
    @{ Cell = 'E32'; Value = '  -1.08%  ' },
    @{ Cell = 'D33'; Value = '1.538' },
    @{ Cell = 'E33'; Value = '  -0.56%  ' },
    @{ Cell = 'E34'; Value = '  +0.02%  ' },
    @{ Cell = 'D35'; Value = '0.8925' },
    @{ Cell = 'E35'; Value = '  -0.19%  ' },
    @{ Cell = 'D36'; Value = '2.579' },
    @{ Cell = 'E36'; Value = '  -1.80%  ' },
    @{ Cell = 'D37'; Value = '1.131.88' },
    @{ Cell = 'E37'; Value = '  -1.72%  ' },
    @{ Cell = 'D38'; Value = '0.5423' },
    @{ Cell = 'E38'; Value = '  -3.11%  ' },
    @{ Cell = 'D39'; Value = '0.01557' },
    @{ Cell = 'E39'; Value = '  -0.70%  ' },
    @{ Cell = 'B40'; Value = 'PaxDollar' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Cell = 'D40'; Value = '1.002' },
    @{ Cell = 'E40'; Value = '  +0.14%  ' },
    @{ Cell = 'B41'; Value = 'mCoin' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin' },
    @{ Cell = 'D41'; Value = '2.545' },
    @{ Cell = 'E41'; Value = '  -0.58%  ' },
    @{ Cell = 'D42'; Value = '0.0₈128' },
    @{ Cell = 'E42'; Value = '  +9.07%  ' },
    @{ Cell = 'D43'; Value = '5.574' },
    @{ Cell = 'E43'; Value = '  -2.45%  ' },
    @{ Cell = 'D44'; Value = '0.8160' },
    @{ Cell = 'E44'; Value = '  +1.22%  ' },
    @{ Cell = 'D45'; Value = '99.25' },
    @{ Cell = 'E45'; Value = '  -0.55%  ' },
    @{ Cell = 'D46'; Value = '1.776.52' },
    @{ Cell = 'E46'; Value = '  -0.23%  ' },
    @{ Cell = 'D47'; Value = '0.4533' },
    @{ Cell = 'E47'; Value = '  +0.54%  ' },
    @{ Cell = 'D48'; Value = '1.002' },
    @{ Cell = 'E48'; Value = '  -0.32%  ' },
    @{ Cell = 'D49'; Value = '54.72' },
    @{ Cell = 'E49'; Value = '  -0.06%  ' },
    @{ Cell = 'E50'; Value = '  +0.60%  ' },
    @{ Cell = 'D51'; Value = '1.005' },
    @{ Cell = 'E51'; Value = '  +0.29%  ' }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Range($u.Cell) $u.Value
}
